$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.663.23'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.597.73'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '211.48'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '0.514'
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D12").Value = '1.822.15'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '1.604.53'
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").Value = '65.16'
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '26.651.27'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '7.01'
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("D24").Value = '8.97'
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("D25").Value = '144.31'
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '7.11'
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("E33").Value = '  +1.58%  '
$ws.Range("D34").Value = '1.286.42'
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("E35").Value = '  -6.50%  '
$ws.Range("D36").Value = '2.45'
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("D39").Value = '0.834'
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  +20.27%  '
$ws.Range("D41").Value = '5.51'
$ws.Range("E41").Value = '  +2.33%  '
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").Value = "'63.60"
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("D45").Value = '1.735.29'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '90.69'
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("E47").Value = '  -3.50%  '
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("D50").Value = '0.0508'
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("E51").Value = '  -0.02%  '
